$d = $word.ActiveDocument

# --- "Solution With MongoDB" (para 6) becomes "Lots of effort" at ilvl 1 ---
$p6 = $d.Paragraphs(6)
$p6.Range.Text = "Lots of effort"
$p6.Range.ListFormat.ListLevelNumber = 2

# --- new paragraph "Not scalable." inserted right after, also ilvl 1 ---
$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs(7)
$p7.Range.Text = "Not scalable."
$p7.Range.ListFormat.ListLevelNumber = 2

# --- "Solution With MongoDB" re-added as its own paragraph at ilvl 0 ---
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
$p8.Range.Text = "Solution With MongoDB"
$p8.Range.ListFormat.ListLevelNumber = 1

# --- "Collections" paragraph (now #10) gains " " + Wingdings arrow symbol + " media alone" ---
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range
$insertPoint10 = $d.Range($r10.End - 1, $r10.End - 1)
$xml10 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> media alone</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint10.InsertXML($xml10) | Out-Null
$endOfPara10 = $r10.End - 1
$markRange10 = $d.Range($endOfPara10, $endOfPara10 + 1)
$markRange10.Delete()

# --- "Data Analytics with MongoDB" (now #12) gains " 5" in its own run ---
$p12 = $d.Paragraphs(12)
$r12 = $p12.Range
$insertPoint12 = $d.Range($r12.End - 1, $r12.End - 1)
$xml12 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> 5</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint12.InsertXML($xml12) | Out-Null
$endOfPara12 = $r12.End - 1
$markRange12 = $d.Range($endOfPara12, $endOfPara12 + 1)
$markRange12.Delete()

# --- two new ilvl-1 bullets appended after "What I learned, challenges." (now #14) ---
$p14 = $d.Paragraphs(14)
$p14.Range.InsertParagraphAfter()
$p15 = $d.Paragraphs(15)
$p15.Range.Text = "Metadata challenges"
$p15.Range.ListFormat.ListLevelNumber = 2

$p15.Range.InsertParagraphAfter()
$p16 = $d.Paragraphs(16)
$p16.Range.Text = "Displaying media in application layer"
$p16.Range.ListFormat.ListLevelNumber = 2
